$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in C3 with "ARREST" (existing shared string) and copy the
# italic/"Courier New" style used by the other cells in that column (e.g. C4).
$ws.Range("C3").Value = "ARREST"
$ws.Range("C4").Copy($ws.Range("C3"))

# Fill in D4 with the new value "в пользу кого 3" (new shared string).
$ws.Range("D4").Value = "в пользу кого 3"

# Update the active selection to D3:D4, with D3 as the active cell.
$ws.Range("D3:D4").Select()
